$wb = $excel.ActiveWorkbook

# --- 1. Update selections on existing sheets that changed in the diff ---
$wsProtein = $wb.Worksheets.Item("Protein")
$wsProtein.Range("G2:H2").Select()

$wsConnM = $wb.Worksheets.Item("Con_nM")
$wsConnM.Range("B2:B9").Select()

$wsPP2C = $wb.Worksheets.Item("PP2C_Activity")
$wsPP2C.Range("B2:M9").Select()

# --- 2. Add the new "Original_Data" sheet at the end of the workbook ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$ws.Name = "Original_Data"

# --- 3. Populate header rows (string order matters for shared-string indices) ---
$ws.Range("B1").Value = "AtPYR1"
$ws.Range("B1").HorizontalAlignment = -4108

$ws.Range("H1").Value = "AtPYL2"
$ws.Range("H1").HorizontalAlignment = -4108

$ws.Range("N1").Value = "Protein"

$ws.Range("B2").Value = "ABA"
$ws.Range("E2").Value = "OP"
$ws.Range("H2").Value = "ABA"
$ws.Range("K2").Value = "OP"
$ws.Range("N2").Value = "Compound"

$ws.Range("B3").Value = 1
$ws.Range("C3").Value = 2
$ws.Range("D3").Value = 3
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 2
$ws.Range("G3").Value = 3
$ws.Range("H3").Value = 1
$ws.Range("I3").Value = 2
$ws.Range("J3").Value = 3
$ws.Range("K3").Value = 1
$ws.Range("L3").Value = 2
$ws.Range("M3").Value = 3

$ws.Range("N3").Value = "Replicate"

$ws.Range("A3").Value = "Con_nM"

# --- 4. Con_nM column (A4:A11), formatted + one shared formula block ---
# Note: alignment is built up one property at a time, always landing on a
# style that already exists in the workbook (2 -> 4 -> 3), to avoid the
# engine minting a throw-away intermediate cellXf that nothing references.
$ws.Range("A11").Value = 0
$ws.Range("A11").HorizontalAlignment = -4108
$ws.Range("A11").VerticalAlignment = -4108

$ws.Range("A4").Value = 10000
$ws.Range("A5").Value = 1000
$ws.Range("A6:A10").Formula = "=A5/3"

$ws.Range("A4:A10").HorizontalAlignment = -4108
$ws.Range("A4:A10").VerticalAlignment = -4108
$ws.Range("A4:A10").WrapText = $true

# --- 5. Data block B4:M11 ---
$data = @(
    @(5.4378080000000004, 5.7129700000000003, 6.1538830000000004, 6.2571149999999998, 8.0371240000000004, 6.3505310000000001, 0.74690500000000004, 3.8373159999999999, 3.6045630000000002, 2.0028009999999998, 3.6759249999999999, 3.6005989999999999),
    @(10.96649, 10.54702, 6.7885369999999998, 8.3268260000000005, 4.4356669999999996, 6.4766620000000001, 12.79105, 11.029299999999999, 10.5486, 1.0182709999999999, 5.7940290000000001, 4.842536),
    @(31.875520000000002, 32.046720000000001, 24.757660000000001, 7.4529959999999997, 6.2952810000000001, 10.06467, 19.1554, 20.9221, 21.753419999999998, 6.913189, 5.5144450000000003, 9.2407120000000003),
    @(65.006529999999998, 63.774299999999997, 51.448360000000001, 17.10474, 15.8278, 19.26642, 38.687019999999997, 38.615160000000003, 39.655439999999999, 6.0091890000000001, 5.7898990000000001, 7.1634520000000004),
    @(82.995639999999995, 81.072779999999995, 76.431010000000001, 36.861510000000003, 40.049309999999998, 42.204810000000002, 66.695350000000005, 61.479480000000002, 68.87585, 18.893160000000002, 19.058350000000001, 24.164359999999999),
    @(95.532409999999999, 90.636200000000002, 91.672139999999999, 63.821550000000002, 70.059039999999996, 74.813479999999998, 83.557090000000002, 82.161230000000003, 84.71754, 52.204079999999998, 55.813470000000002, 58.75385),
    @(98.207689999999999, 86.423349999999999, 95.354299999999995, 92.599040000000002, 90.672550000000001, 90.650739999999999, 92.464939999999999, 80.521720000000002, 88.331069999999997, 84.556479999999993, 80.942959999999999, 83.552959999999999),
    @(102.16249999999999, 113.4924, 86.48151, 89.858329999999995, 100.8139, 93.031599999999997, 99.241849999999999, 97.792310000000001, 100.0761, 106.33669999999999, 98.655420000000007, 93.757549999999995)
)

$cols = @("B", "C", "D", "E", "F", "G", "H", "I", "J", "K", "L", "M")
for ($r = 0; $r -lt $data.Length; $r++) {
    $rowNum = 4 + $r
    $rowVals = $data[$r]
    for ($c = 0; $c -lt $cols.Length; $c++) {
        $ws.Range($cols[$c] + $rowNum).Value = $rowVals[$c]
    }
}

# Build up to style 5 (wrap + vertical-center, general horizontal) via
# styles 2 -> 4 -> 3 (all already present) and finish by clearing the
# horizontal override back to general, landing exactly on style 5.
$ws.Range("B4:M11").HorizontalAlignment = -4108
$ws.Range("B4:M11").VerticalAlignment = -4108
$ws.Range("B4:M11").WrapText = $true
$ws.Range("B4:M11").HorizontalAlignment = 1

# --- 6. Final selection / active sheet state ---
$ws.Range("I23").Select()

Write-Output "done"
